$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-02 Sunday" "2024-06-03 Monday"

Replace-Text "185×7=" "766×5="
Replace-Text "829×8=" "606×7="
Replace-Text "391×4=" "894×8="
Replace-Text "973×7=" "531×8="
Replace-Text "832×9=" "603×8="

Replace-Text "177×7=" "111×9="
Replace-Text "903×9=" "889×2="
Replace-Text "408×2=" "286×9="
Replace-Text "466×5=" "872×2="
Replace-Text "580×2=" "514×7="

Replace-Text "982×3=" "219×8="
Replace-Text "285×5=" "493×3="
Replace-Text "602×7=" "310×9="
Replace-Text "927×5=" "317×8="
Replace-Text "748×9=" "832×3="

Replace-Text "383×2=" "362×5="
Replace-Text "386×8=" "420×2="
Replace-Text "563×5=" "983×8="
Replace-Text "245×2=" "944×6="
Replace-Text "792×4=" "785×3="

Replace-Text "675×7=" "844×7="
Replace-Text "211×7=" "323×3="
Replace-Text "289×7=" "914×8="
Replace-Text "542×2=" "683×8="
Replace-Text "374×4=" "685×9="
